$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (division problems)
$t.Cell(1, 1).Range.Text = "62÷2=31, 0"
$t.Cell(1, 2).Range.Text = "30÷5=6, 0"
$t.Cell(1, 3).Range.Text = "33÷5=6, 3"
$t.Cell(1, 4).Range.Text = "80÷3=26, 2"
$t.Cell(1, 5).Range.Text = "66÷8=8, 2"

# Row 5 (division problems)
$t.Cell(5, 1).Range.Text = "81÷8=10, 1"
$t.Cell(5, 2).Range.Text = "23÷2=11, 1"
$t.Cell(5, 3).Range.Text = "89÷3=29, 2"
$t.Cell(5, 4).Range.Text = "72÷8=9, 0"
$t.Cell(5, 5).Range.Text = "48÷5=9, 3"

# Row 9 (division problems)
$t.Cell(9, 1).Range.Text = "94÷5=18, 4"
$t.Cell(9, 2).Range.Text = "12÷8=1, 4"
$t.Cell(9, 3).Range.Text = "91÷3=30, 1"
$t.Cell(9, 4).Range.Text = "90÷2=45, 0"
$t.Cell(9, 5).Range.Text = "71÷9=7, 8"

# Row 13 (division problems)
$t.Cell(13, 1).Range.Text = "16÷8=2, 0"
$t.Cell(13, 2).Range.Text = "52÷4=13, 0"
$t.Cell(13, 3).Range.Text = "10÷9=1, 1"
$t.Cell(13, 4).Range.Text = "37÷8=4, 5"
$t.Cell(13, 5).Range.Text = "17÷2=8, 1"

# Row 17 (division problems)
$t.Cell(17, 1).Range.Text = "26÷9=2, 8"
$t.Cell(17, 2).Range.Text = "45÷7=6, 3"
$t.Cell(17, 3).Range.Text = "20÷2=10, 0"
$t.Cell(17, 4).Range.Text = "32÷3=10, 2"
$t.Cell(17, 5).Range.Text = "52÷2=26, 0"
